$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3D")

$ws.Range("A6").Value = "2026-02-09 18:26:25"
$ws.Range("B6").Value = "Alamin hamza"
$ws.Range("C6").Value = "Number 8"
$ws.Range("D6").Value = 7
